# feat: update league_data.xlsx with new league statistics
#
# "champs" sheet (sheet5.xml):
#   - years 2016-2007 (rows 11-20): venue corrected from "Sacred Heart" to
#     "Magna Centre"
#   - nine more seasons appended (years 2005 down to 1997, rows 22-30),
#     continuing the Team 2..Team 6 rotation, all played at "Sacred Heart"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("champs")
$ws.Activate()

# Correct the venue for 2016 .. 2007 (rows 11-20): Sacred Heart -> Magna Centre
for ($r = 11; $r -le 20; $r++) {
    $ws.Cells.Item($r, 3).Value = "Magna Centre"
}

# Append historical seasons 2005 .. 1997 (rows 22-30)
$teams = @("Team 2", "Team 3", "Team 4", "Team 5", "Team 6")

$row = 22
$year = 2005
for ($i = 0; $i -lt 9; $i++) {
    $ws.Cells.Item($row, 1).Value = $year
    $ws.Cells.Item($row, 2).Value = $teams[$i % 5]
    $ws.Cells.Item($row, 3).Value = "Sacred Heart"
    $row++
    $year--
}

# Widen column D slightly to keep fitting the player-name entries
$ws.Columns.Item(4).ColumnWidth = 10.6

# Leave the selection on the next empty row beneath the appended table
$ws.Range("B31").Select() | Out-Null

# Restore portrait page setup for the sheet (picked up when the view/print
# area was touched)
$ws.PageSetup.Orientation = 1
